# RegistrationTestdata.xlsx update
# - Refresh the sample rows on both the "Positve Testdata" and
#   "Negative Testdata" sheets with new fullname / email / password values.
# - "Negative Testdata" row 2's password/confirm columns switch from numeric
#   "123" to the text value "sk123".
# - Auto-fit the contactno column on "Negative Testdata" (it now holds a
#   shorter 3-digit sample value next to a 10-digit one).
# - The sheet that is active/selected flips from "Negative Testdata" to
#   "Positve Testdata", and the last selection on each sheet becomes the
#   (now empty) row below the data, A4:XFD4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Positve Testdata")
$ws2 = $wb.Worksheets.Item("Negative Testdata")

# --- Positve Testdata ---------------------------------------------------
$ws1.Range("A2").Value = "fake90"
$ws1.Range("B2").Value = "fake879@g.com"
$ws1.Range("C2").Value = 1267438191
$ws1.Range("D2").Value = "fake@12345"
$ws1.Range("E2").Value = "fake@12345"

$ws1.Range("A3").Value = "fake91"
$ws1.Range("B3").Value = "fake456@g.com"
$ws1.Range("C3").Value = 3456789012
$ws1.Range("D3").Value = "fake@12345"
$ws1.Range("E3").Value = "fake@12345"

# --- Negative Testdata ----------------------------------------------------
$ws2.Range("A2").Value = "fake12"
$ws2.Range("B2").Value = "sk@g.com"
$ws2.Range("C2").Value = 123
$ws2.Range("D2").Value = "sk123"
$ws2.Range("E2").Value = "sk123"

$ws2.Range("A3").Value = "fake2"
$ws2.Range("B3").Value = "new890@g.com"
$ws2.Range("C3").Value = 1234567890
$ws2.Range("D3").Value = "new@1234"
$ws2.Range("E3").Value = "new@123"

# Auto-fit the contactno column now that it holds mixed-length numbers
# (10 characters wide to fit "1234567890").
$ws2.Columns.Item(3).ColumnWidth = 10

# Selection / active-sheet bookkeeping: the workbook now opens on
# "Positve Testdata" (tab 1) with the row below the table selected on both
# sheets, instead of "Negative Testdata" being the active tab with a stray
# Q10 selection.
[void]$ws2.Activate()
[void]$ws2.Range("A4:XFD4").Select()

[void]$ws1.Activate()
[void]$ws1.Range("A4:XFD4").Select()
